# Update stats for 2025-11 (row 24 in Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B24").Value = 6400
$ws.Range("C24").Value = 1004
$ws.Range("D24").Value = 5975833
$ws.Range("E24").Value = 933.72390625
$ws.Range("F24").Value = 9.103307193999321
$ws.Range("G24").Value = 4.041450777202082
$ws.Range("H24").Value = 26.59129020805107
